# Weekly update: insert a new price record as row 14, pushing the
# existing rows 14-29 down to 15-30 (new dimension A1:T30).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 14 - this shifts the existing
# rows 14..29 down to 15..30, preserving all their data/styles.
$ws.Rows("14:14").Insert()

# Populate the newly inserted row 14 with the new weekly record.
$ws.Range("A14").Value = 4
$ws.Range("B14").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C14").Value = "Los Lagos"
$ws.Range("D14").Value = 44904
$ws.Range("E14").Value = 10
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100101
$ws.Range("H14").Value = "Berries"
$ws.Range("I14").Value = 100101001
$ws.Range("J14").Value = "Arándano (blue)"
$ws.Range("K14").Value = "Sin especificar"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 400
$ws.Range("N14").Value = 5000
$ws.Range("O14").Value = 5500
$ws.Range("P14").Value = 5250
$ws.Range("Q14").Value = "$/bandeja 2 kilos"
$ws.Range("R14").Value = "Provincia de Curicó"
$ws.Range("S14").Value = 2625
$ws.Range("T14").Value = 2
